$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.657.88"
$ws.Range("E2").Value = "  -4.43%  "

# Row 3
$ws.Range("D3").Value = "2.966.59"
$ws.Range("E3").Value = "  -6.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.60"
$ws.Range("E5").Value = "  -5.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.45"
$ws.Range("E6").Value = "  -7.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -3.61%  "

# Row 9
$ws.Range("D9").Value = "2.974.13"
$ws.Range("E9").Value = "  -6.15%  "

# Row 10
$ws.Range("E10").Value = "  -4.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.10"
$ws.Range("E11").Value = "  -8.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.366"
$ws.Range("E12").Value = "  -4.81%  "

# Row 13
$ws.Range("D13").Value = "3.485.43"
$ws.Range("E13").Value = "  -6.49%  "

# Row 14
$ws.Range("E14").Value = "  -3.33%  "

# Row 15
$ws.Range("D15").Value = "61.679.37"
$ws.Range("E15").Value = "  -4.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.67"

# Row 17
$ws.Range("D17").Value = "2.970.72"
$ws.Range("E17").Value = "  -5.96%  "

# Row 18
$ws.Range("E18").Value = "  -5.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.14"
$ws.Range("E19").Value = "  -2.50%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.97"
$ws.Range("E20").Value = "  -6.21%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.36"
$ws.Range("E21").Value = "  -6.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("E22").Value = "  -6.35%  "

# Row 23
$ws.Range("E23").Value = "  -0.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.10"
$ws.Range("E24").Value = "  -5.47%  "

# Row 25
$ws.Range("E25").Value = "  -3.39%  "

# Row 26
$ws.Range("D26").Value = "3.093.72"
$ws.Range("E26").Value = "  -6.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.186"
$ws.Range("E27").Value = "  -4.44%  "

# Row 28
$ws.Range("E28").Value = "  -0.19%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0929"
$ws.Range("E29").Value = "  -9.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.30"
$ws.Range("E30").Value = "  -6.43%  "

# Row 31
$ws.Range("E31").Value = "  -0.05%  "

# Row 32
$ws.Range("E32").Value = "  -5.85%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.40"
$ws.Range("E33").Value = "  -4.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.75"
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
$ws.Range("E35").Value = "  -5.25%  "

# Row 36
$ws.Range("E36").Value = "  -6.67%  "

# Row 37
$ws.Range("E37").Value = "  -5.58%  "

# Row 38
$ws.Range("E38").Value = "  -5.52%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  -9.12%  "

# Row 40
$ws.Range("E40").Value = "  -4.75%  "

# Row 41
$ws.Range("D41").Value = "2.414.29"
$ws.Range("E41").Value = "  -10.25%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.07"
$ws.Range("E42").Value = "  -4.00%  "

# Row 43
$ws.Range("E43").Value = "  -7.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.664"
$ws.Range("E44").Value = "  -4.72%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("E45").Value = "  -4.72%  "

# Row 46
$ws.Range("E46").Value = "  -0.32%  "

# Row 47
$ws.Range("E47").Value = "  -5.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("E48").Value = "  -9.77%  "

# Row 49
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "268.03"
$ws.Range("E49").Value = "  -8.01%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0950"
$ws.Range("E50").Value = "  -3.65%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.67"
$ws.Range("E51").Value = "  -8.63%  "
